# accelpressure done. work on making figure
# Update foundPath data points: several existing rows' Y (B), Velocity (C), and
# Cost (D) values are recalculated/updated, and a new data row (72) is appended
# to the path-node table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, "B").Value = 690
$ws.Cells.Item(4, "C").Value = 19
$ws.Cells.Item(4, "D").Value = 5.925925925925926
$ws.Cells.Item(5, "B").Value = 700
$ws.Cells.Item(5, "C").Value = 21
$ws.Cells.Item(5, "D").Value = 6.425925925925926
$ws.Cells.Item(6, "B").Value = 740
$ws.Cells.Item(6, "C").Value = 26
$ws.Cells.Item(6, "D").Value = 8.128053585500394
$ws.Cells.Item(7, "A").Value = 900
$ws.Cells.Item(7, "B").Value = 810
$ws.Cells.Item(7, "C").Value = 33
$ws.Cells.Item(7, "D").Value = 10.52502572511581
$ws.Cells.Item(8, "A").Value = 900
$ws.Cells.Item(8, "B").Value = 870
$ws.Cells.Item(8, "C").Value = 32
$ws.Cells.Item(8, "D").Value = 12.37117957126966
$ws.Cells.Item(9, "A").Value = 900
$ws.Cells.Item(9, "B").Value = 930
$ws.Cells.Item(9, "C").Value = 30
$ws.Cells.Item(9, "D").Value = 14.3066634422374
$ws.Cells.Item(10, "B").Value = 970
$ws.Cells.Item(10, "C").Value = 27
$ws.Cells.Item(10, "D").Value = 15.7533671705243
$ws.Cells.Item(11, "B").Value = 1020
$ws.Cells.Item(11, "C").Value = 21
$ws.Cells.Item(11, "D").Value = 17.99718584016367
$ws.Cells.Item(12, "B").Value = 1030
$ws.Cells.Item(12, "C").Value = 17
$ws.Cells.Item(12, "D").Value = 18.74150876772846
$ws.Cells.Item(13, "A").Value = 820
$ws.Cells.Item(13, "B").Value = 1050
$ws.Cells.Item(13, "C").Value = 11
$ws.Cells.Item(13, "D").Value = 21.93589159272816
$ws.Cells.Item(14, "B").Value = 1050
$ws.Cells.Item(14, "C").Value = 13
$ws.Cells.Item(14, "D").Value = 22.76922492606149
$ws.Cells.Item(15, "A").Value = 770
$ws.Cells.Item(15, "B").Value = 1030
$ws.Cells.Item(15, "C").Value = 17
$ws.Cells.Item(15, "D").Value = 25.75064889606121
$ws.Cells.Item(16, "A").Value = 750
$ws.Cells.Item(16, "B").Value = 1010
$ws.Cells.Item(16, "C").Value = 21
$ws.Cells.Item(16, "D").Value = 27.23929475119079
$ws.Cells.Item(17, "A").Value = 730
$ws.Cells.Item(17, "B").Value = 970
$ws.Cells.Item(17, "C").Value = 21
$ws.Cells.Item(17, "D").Value = 29.36888330119059
$ws.Cells.Item(18, "A").Value = 730
$ws.Cells.Item(18, "B").Value = 960
$ws.Cells.Item(18, "C").Value = 17
$ws.Cells.Item(18, "D").Value = 29.89519909066427
$ws.Cells.Item(19, "A").Value = 730
$ws.Cells.Item(19, "B").Value = 880
$ws.Cells.Item(19, "C").Value = 25
$ws.Cells.Item(19, "D").Value = 33.70472290018808
$ws.Cells.Item(20, "A").Value = 720
$ws.Cells.Item(20, "B").Value = 860
$ws.Cells.Item(20, "C").Value = 17
$ws.Cells.Item(20, "D").Value = 34.76951717518798
$ws.Cells.Item(21, "A").Value = 700
$ws.Cells.Item(21, "B").Value = 850
$ws.Cells.Item(21, "C").Value = 8
$ws.Cells.Item(21, "D").Value = 36.55837155718781
$ws.Cells.Item(22, "A").Value = 680
$ws.Cells.Item(22, "B").Value = 850
$ws.Cells.Item(22, "C").Value = 14
$ws.Cells.Item(22, "D").Value = 38.37655337536963
$ws.Cells.Item(23, "A").Value = 670
$ws.Cells.Item(23, "B").Value = 850
$ws.Cells.Item(23, "C").Value = 16
$ws.Cells.Item(23, "D").Value = 39.0432200420363
$ws.Cells.Item(24, "A").Value = 650
$ws.Cells.Item(24, "B").Value = 850
$ws.Cells.Item(24, "C").Value = 19
$ws.Cells.Item(24, "D").Value = 40.18607718489344
$ws.Cells.Item(25, "A").Value = 640
$ws.Cells.Item(25, "B").Value = 850
$ws.Cells.Item(25, "C").Value = 21
$ws.Cells.Item(25, "D").Value = 40.68607718489344
$ws.Cells.Item(26, "A").Value = 620
$ws.Cells.Item(26, "B").Value = 850
$ws.Cells.Item(26, "C").Value = 24
$ws.Cells.Item(26, "D").Value = 41.57496607378233
$ws.Cells.Item(27, "A").Value = 570
$ws.Cells.Item(27, "C").Value = 21
$ws.Cells.Item(27, "D").Value = 43.968372654731
$ws.Cells.Item(28, "A").Value = 490
$ws.Cells.Item(28, "B").Value = 890
$ws.Cells.Item(28, "C").Value = 29
$ws.Cells.Item(28, "D").Value = 47.26685715522513
$ws.Cells.Item(29, "A").Value = 460
$ws.Cells.Item(29, "B").Value = 890
$ws.Cells.Item(29, "C").Value = 17
$ws.Cells.Item(29, "D").Value = 48.57120498131208
$ws.Cells.Item(30, "B").Value = 860
$ws.Cells.Item(30, "C").Value = 17
$ws.Cells.Item(30, "D").Value = 52.517207294547
$ws.Cells.Item(31, "A").Value = 380
$ws.Cells.Item(31, "B").Value = 840
$ws.Cells.Item(31, "C").Value = 21
$ws.Cells.Item(31, "D").Value = 54.00585314967658
$ws.Cells.Item(32, "A").Value = 360
$ws.Cells.Item(32, "B").Value = 770
$ws.Cells.Item(32, "C").Value = 15.17368524526475
$ws.Cells.Item(32, "D").Value = 58.03093925296958
$ws.Cells.Item(33, "A").Value = 360
$ws.Cells.Item(33, "B").Value = 760
$ws.Cells.Item(33, "C").Value = 16
$ws.Cells.Item(33, "D").Value = 58.67250600526872
$ws.Cells.Item(34, "A").Value = 370
$ws.Cells.Item(34, "B").Value = 720
$ws.Cells.Item(34, "C").Value = 22
$ws.Cells.Item(34, "D").Value = 60.84256159769906
$ws.Cells.Item(35, "A").Value = 390
$ws.Cells.Item(35, "B").Value = 690
$ws.Cells.Item(35, "C").Value = 26
$ws.Cells.Item(35, "D").Value = 62.34487462914239
$ws.Cells.Item(36, "A").Value = 440
$ws.Cells.Item(36, "B").Value = 630
$ws.Cells.Item(36, "C").Value = 34
$ws.Cells.Item(36, "D").Value = 64.94829118777794
$ws.Cells.Item(37, "A").Value = 510
$ws.Cells.Item(37, "B").Value = 550
$ws.Cells.Item(37, "C").Value = 30
$ws.Cells.Item(37, "D").Value = 68.27021175425752
$ws.Cells.Item(38, "A").Value = 540
$ws.Cells.Item(38, "B").Value = 500
$ws.Cells.Item(38, "C").Value = 23
$ws.Cells.Item(38, "D").Value = 70.47057095985951
$ws.Cells.Item(39, "A").Value = 550
$ws.Cells.Item(39, "B").Value = 460
$ws.Cells.Item(39, "C").Value = 19
$ws.Cells.Item(39, "D").Value = 72.43395459110602
$ws.Cells.Item(40, "A").Value = 550
$ws.Cells.Item(40, "B").Value = 450
$ws.Cells.Item(40, "C").Value = 21
$ws.Cells.Item(40, "D").Value = 72.93395459110602
$ws.Cells.Item(41, "A").Value = 540
$ws.Cells.Item(41, "B").Value = 440
$ws.Cells.Item(41, "C").Value = 1
$ws.Cells.Item(41, "D").Value = 74.21960328417246
$ws.Cells.Item(42, "A").Value = 520
$ws.Cells.Item(42, "B").Value = 430
$ws.Cells.Item(42, "C").Value = 12
$ws.Cells.Item(42, "D").Value = 77.65970786494137
$ws.Cells.Item(43, "A").Value = 500
$ws.Cells.Item(43, "B").Value = 430
$ws.Cells.Item(43, "C").Value = 16
$ws.Cells.Item(43, "D").Value = 79.0882792935128
$ws.Cells.Item(44, "A").Value = 490
$ws.Cells.Item(44, "B").Value = 430
$ws.Cells.Item(44, "C").Value = 18
$ws.Cells.Item(44, "D").Value = 79.67651458763045
$ws.Cells.Item(45, "A").Value = 480
$ws.Cells.Item(45, "B").Value = 430
$ws.Cells.Item(45, "C").Value = 20
$ws.Cells.Item(45, "D").Value = 80.20283037710414
$ws.Cells.Item(46, "A").Value = 460
$ws.Cells.Item(46, "B").Value = 430
$ws.Cells.Item(46, "C").Value = 23
$ws.Cells.Item(46, "D").Value = 81.13306293524367
$ws.Cells.Item(47, "A").Value = 440
$ws.Cells.Item(47, "B").Value = 430
$ws.Cells.Item(47, "C").Value = 25
$ws.Cells.Item(47, "D").Value = 81.966396268577
$ws.Cells.Item(48, "A").Value = 400
$ws.Cells.Item(48, "B").Value = 440
$ws.Cells.Item(48, "C").Value = 24
$ws.Cells.Item(48, "D").Value = 83.64929652393114
$ws.Cells.Item(49, "A").Value = 350
$ws.Cells.Item(49, "B").Value = 430
$ws.Cells.Item(49, "C").Value = 18.47006582938073
$ws.Cells.Item(49, "D").Value = 86.05052638501851
$ws.Cells.Item(50, "A").Value = 290
$ws.Cells.Item(50, "B").Value = 390
$ws.Cells.Item(50, "C").Value = 21
$ws.Cells.Item(50, "D").Value = 89.7044866201495
$ws.Cells.Item(51, "A").Value = 280
$ws.Cells.Item(51, "B").Value = 380
$ws.Cells.Item(51, "C").Value = 23
$ws.Cells.Item(51, "D").Value = 90.34731096668273
$ws.Cells.Item(52, "A").Value = 260
$ws.Cells.Item(52, "B").Value = 350
$ws.Cells.Item(52, "C").Value = 22
$ws.Cells.Item(52, "D").Value = 91.94977820022228
$ws.Cells.Item(53, "A").Value = 240
$ws.Cells.Item(53, "B").Value = 310
$ws.Cells.Item(53, "C").Value = 27
$ws.Cells.Item(53, "D").Value = 93.77513981450782
$ws.Cells.Item(54, "A").Value = 240
$ws.Cells.Item(54, "B").Value = 290
$ws.Cells.Item(54, "C").Value = 17
$ws.Cells.Item(54, "D").Value = 94.68423072359873
$ws.Cells.Item(55, "A").Value = 260
$ws.Cells.Item(55, "B").Value = 260
$ws.Cells.Item(55, "C").Value = 11
$ws.Cells.Item(55, "D").Value = 97.25962449178729
$ws.Cells.Item(56, "A").Value = 270
$ws.Cells.Item(56, "B").Value = 240
$ws.Cells.Item(56, "C").Value = 16
$ws.Cells.Item(56, "D").Value = 98.91597114178714
$ws.Cells.Item(57, "A").Value = 290
$ws.Cells.Item(57, "B").Value = 220
$ws.Cells.Item(57, "C").Value = 21
$ws.Cells.Item(57, "D").Value = 100.444850668677
$ws.Cells.Item(58, "A").Value = 330
$ws.Cells.Item(58, "B").Value = 170
$ws.Cells.Item(58, "D").Value = 103.0061003636501
$ws.Cells.Item(59, "A").Value = 360
$ws.Cells.Item(59, "B").Value = 140
$ws.Cells.Item(59, "C").Value = 33
$ws.Cells.Item(59, "D").Value = 104.3746941336886
$ws.Cells.Item(60, "A").Value = 380
$ws.Cells.Item(60, "B").Value = 130
$ws.Cells.Item(60, "D").Value = 105.1457520569644
$ws.Cells.Item(61, "A").Value = 480
$ws.Cells.Item(61, "B").Value = 110
$ws.Cells.Item(61, "C").Value = 27.43989448617086
$ws.Cells.Item(61, "D").Value = 109.0351721738851
$ws.Cells.Item(62, "A").Value = 500
$ws.Cells.Item(62, "B").Value = 110
$ws.Cells.Item(62, "C").Value = 30
$ws.Cells.Item(62, "D").Value = 109.7315522831771
$ws.Cells.Item(68, "A").Value = 900
$ws.Cells.Item(68, "B").Value = 290
$ws.Cells.Item(68, "C").Value = 27
$ws.Cells.Item(68, "D").Value = 126.4166584414511
$ws.Cells.Item(69, "A").Value = 920
$ws.Cells.Item(69, "B").Value = 360
$ws.Cells.Item(69, "C").Value = 34
$ws.Cells.Item(69, "D").Value = 128.8035797166251
$ws.Cells.Item(70, "A").Value = 940
$ws.Cells.Item(70, "C").Value = 29
$ws.Cells.Item(70, "D").Value = 132.6656511644335
$ws.Cells.Item(71, "A").Value = 940
$ws.Cells.Item(71, "B").Value = 500
$ws.Cells.Item(71, "C").Value = 31
$ws.Cells.Item(71, "D").Value = 133.3323178311001
$ws.Cells.Item(72, "A").Value = 940
$ws.Cells.Item(72, "B").Value = 530
$ws.Cells.Item(72, "C").Value = 34
$ws.Cells.Item(72, "D").Value = 134.2553947541771

